$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 115.4136146666667
$ws.Range("H2").Value = 346.240844
$ws.Range("I2").Value = 0.2619217538490851
$ws.Range("J2").Value = 0.2619217538490851
$ws.Range("M2").Value = 2.425633666666667
$ws.Range("N2").Value = 7.276901000000001
$ws.Range("O2").Value = 0.0662600404061536
$ws.Range("P2").Value = 0.06626004040615362
$ws.Range("Q2").Value = 279.9511493271604
$ws.Range("R2").Value = 2519.560343944444
$ws.Range("S2").Value = 0.01735494599329099
$ws.Range("T2").Value = 0.01735494599329099
$ws.Range("G3").Value = 115.4136146666667
$ws.Range("H3").Value = 346.240844
$ws.Range("I3").Value = 0.2619217538490851
$ws.Range("J3").Value = 0.2619217538490851
$ws.Range("O3").Value = 0.4234968256437875
$ws.Range("P3").Value = 0.4234968256437876
$ws.Range("Q3").Value = 1789.289930230284
$ws.Range("R3").Value = 16103.60937207256
$ws.Range("S3").Value = 0.110923031322141
$ws.Range("T3").Value = 0.1109230313221411
$ws.Range("G4").Value = 115.4136146666667
$ws.Range("H4").Value = 346.240844
$ws.Range("I4").Value = 0.2619217538490851
$ws.Range("J4").Value = 0.2619217538490851
$ws.Range("M4").Value = 18.67887366666667
$ws.Range("N4").Value = 56.03662100000001
$ws.Range("O4").Value = 0.5102431339500588
$ws.Range("P4").Value = 0.5102431339500588
$ws.Range("Q4").Value = 2155.796327772014
$ws.Range("R4").Value = 19402.16694994812
$ws.Range("S4").Value = 0.133643776533653
$ws.Range("T4").Value = 0.1336437765336531
$ws.Range("I5").Value = 0.6414314537852458
$ws.Range("J5").Value = 0.6414314537852458
$ws.Range("M5").Value = 2.425633666666667
$ws.Range("N5").Value = 7.276901000000001
$ws.Range("O5").Value = 0.0662600404061536
$ws.Range("P5").Value = 0.06626004040615362
$ws.Range("Q5").Value = 685.5844161964336
$ws.Range("R5").Value = 6170.259745767902
$ws.Range("S5").Value = 0.04250127404558823
$ws.Range("T5").Value = 0.04250127404558824
$ws.Range("I6").Value = 0.6414314537852458
$ws.Range("J6").Value = 0.6414314537852458
$ws.Range("O6").Value = 0.4234968256437875
$ws.Range("P6").Value = 0.4234968256437876
$ws.Range("S6").Value = 0.2716441845461314
$ws.Range("T6").Value = 0.2716441845461314
$ws.Range("I7").Value = 0.6414314537852458
$ws.Range("J7").Value = 0.6414314537852458
$ws.Range("M7").Value = 18.67887366666667
$ws.Range("N7").Value = 56.03662100000001
$ws.Range("O7").Value = 0.5102431339500588
$ws.Range("P7").Value = 0.5102431339500588
$ws.Range("Q7").Value = 5279.422393393262
$ws.Range("R7").Value = 47514.80154053935
$ws.Range("S7").Value = 0.3272859951935261
$ws.Range("T7").Value = 0.3272859951935261
$ws.Range("I8").Value = 0.09664679236566912
$ws.Range("J8").Value = 0.09664679236566913
$ws.Range("M8").Value = 2.425633666666667
$ws.Range("N8").Value = 7.276901000000001
$ws.Range("O8").Value = 0.0662600404061536
$ws.Range("P8").Value = 0.06626004040615362
$ws.Range("Q8").Value = 103.2994785806983
$ws.Range("R8").Value = 929.695307226285
$ws.Range("S8").Value = 0.006403820367274373
$ws.Range("T8").Value = 0.006403820367274375
$ws.Range("I9").Value = 0.09664679236566912
$ws.Range("J9").Value = 0.09664679236566913
$ws.Range("O9").Value = 0.4234968256437875
$ws.Range("P9").Value = 0.4234968256437876
$ws.Range("S9").Value = 0.04092960977551511
$ws.Range("T9").Value = 0.04092960977551512
$ws.Range("I10").Value = 0.09664679236566912
$ws.Range("J10").Value = 0.09664679236566913
$ws.Range("M10").Value = 18.67887366666667
$ws.Range("N10").Value = 56.03662100000001
$ws.Range("O10").Value = 0.5102431339500588
$ws.Range("P10").Value = 0.5102431339500588
$ws.Range("Q10").Value = 795.4696278984984
$ws.Range("R10").Value = 7159.226651086486
$ws.Range("S10").Value = 0.04931336222287963
$ws.Range("T10").Value = 0.04931336222287964
